$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("B4").Value = "MCT-2A-Programação de Computadores"
$ws.Range("F4").Value = "-"

$ws.Range("B6").Value = "MCT-2A-Programação de Computadores"

$ws.Range("B7").Value = "MCT-2A-Programação de Computadores"
$ws.Range("D7").Value = "MCT-2A-Circuitos Elétricos 2"

$ws.Range("B8").Value = "MCT-2A-Programação de Computadores"
$ws.Range("D8").Value = "MCT-2A-Circuitos Elétricos 2"
